# ---------------------------------------------------------------------------
# Apply the edits described by the commit "Add files via upload":
#   1. Move the "Sheet1" tab to the end of the workbook (after "court").
#   2. On "base_parameter": mark A24 with the participation circle, and
#      append six new member rows (92-97).
#   3. On "Sheet1" (the level/court helper sheet): add a small gender-count
#      summary table in columns L:N, rows 4-7.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# 1. Reorder sheets: member, base_parameter, history, court, Sheet1
$sheet1 = $wb.Worksheets.Item("Sheet1")
$court  = $wb.Worksheets.Item("court")
$sheet1.Move($null, $court)

# 2. base_parameter updates ------------------------------------------------
$bp = $wb.Worksheets.Item("base_parameter")

# Mark A24 as participating
$bp.Range("A24").Value = "〇"

# New rows 92-97
$bp.Range("B92").Value = "[FREE_3]"
$bp.Range("C92").Value = "男"
$bp.Range("D92").Value = 3
$bp.Range("F92").Value = "〇"

$bp.Range("B93").Value = "[FREE_3]"
$bp.Range("C93").Value = "男"
$bp.Range("D93").Value = 3
$bp.Range("F93").Value = "〇"

$bp.Range("A94").Value = "〇"
$bp.Range("B94").Value = "[FREE_3]"
$bp.Range("C94").Value = "女"
$bp.Range("D94").Value = 3
$bp.Range("F94").Value = "〇"

$bp.Range("A95").Value = "〇"
$bp.Range("B95").Value = "[FREE_3]"
$bp.Range("C95").Value = "女"
$bp.Range("D95").Value = 3
$bp.Range("F95").Value = "〇"

$bp.Range("A96").Value = "〇"
$bp.Range("B96").Value = "[FREE_3]"
$bp.Range("C96").Value = "女"
$bp.Range("D96").Value = 3
$bp.Range("F96").Value = "〇"

$bp.Range("B97").Value = "[FREE_3]"
$bp.Range("C97").Value = "男"
$bp.Range("D97").Value = 3
$bp.Range("F97").Value = "〇"

# Give the new cells the same centred style as the rest of the table
# (only columns A, D, E, F carry the centred style; B/C stay default,
# matching the existing rows above them)
$bp.Range("A92:A97").HorizontalAlignment = -4108
$bp.Range("A92:A97").VerticalAlignment = -4108
$bp.Range("D92:F97").HorizontalAlignment = -4108
$bp.Range("D92:F97").VerticalAlignment = -4108

# 3. Sheet1 (level chart) gender-count summary table ------------------------
$s1 = $wb.Worksheets.Item("Sheet1")

$s1.Range("M4").Value = "男"
$s1.Range("N4").Value = "女"

$s1.Range("L5").Value = "3.5以上"
$s1.Range("M5").Value = 3
$s1.Range("N5").Value = 4

$s1.Range("L6").Value = "2.5～3.5"
$s1.Range("M6").Value = 0
$s1.Range("N6").Value = 1

$s1.Range("L7").Value = "2.4以下"
$s1.Range("M7").Value = 1
$s1.Range("N7").Value = 0

# Leave the selection on the newly-edited summary table
$sel = $excel.Union($s1.Range("M7"), $s1.Range("N6"), $s1.Range("M5"))
$sel.Select()

# base_parameter ends up with the new last row selected
$bp.Activate()
$bp.Range("A92").Select()

Write-Output "done"
